$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Apply updated AgTests (F) / AgPosit (G) values for rows 303-390
$ws.Cells.Item(303, 6).Value = 9661
$ws.Cells.Item(303, 7).Value = 614
$ws.Cells.Item(304, 6).Value = 6100
$ws.Cells.Item(304, 7).Value = 526
$ws.Cells.Item(305, 6).Value = 3363
$ws.Cells.Item(305, 7).Value = 272
$ws.Cells.Item(306, 6).Value = 75354
$ws.Cells.Item(306, 7).Value = 7667
$ws.Cells.Item(307, 6).Value = 75380
$ws.Cells.Item(307, 7).Value = 6329
$ws.Cells.Item(308, 6).Value = 15552
$ws.Cells.Item(308, 7).Value = 1057
$ws.Cells.Item(309, 6).Value = 77892
$ws.Cells.Item(310, 6).Value = 79239
$ws.Cells.Item(310, 7).Value = 4063
$ws.Cells.Item(311, 6).Value = 61493
$ws.Cells.Item(312, 6).Value = 28172
$ws.Cells.Item(312, 7).Value = 925
$ws.Cells.Item(313, 6).Value = 75804
$ws.Cells.Item(313, 7).Value = 3459
$ws.Cells.Item(314, 6).Value = 64243
$ws.Cells.Item(314, 7).Value = 3144
$ws.Cells.Item(315, 6).Value = 56384
$ws.Cells.Item(315, 7).Value = 2630
$ws.Cells.Item(316, 6).Value = 50711
$ws.Cells.Item(316, 7).Value = 2298
$ws.Cells.Item(317, 6).Value = 63738
$ws.Cells.Item(317, 7).Value = 2175
$ws.Cells.Item(318, 6).Value = 49337
$ws.Cells.Item(318, 7).Value = 1132
$ws.Cells.Item(319, 6).Value = 41358
$ws.Cells.Item(319, 7).Value = 1634
$ws.Cells.Item(320, 6).Value = 71846
$ws.Cells.Item(320, 7).Value = 3306
$ws.Cells.Item(321, 6).Value = 89598
$ws.Cells.Item(321, 7).Value = 2665
$ws.Cells.Item(322, 6).Value = 109364
$ws.Cells.Item(322, 7).Value = 2328
$ws.Cells.Item(323, 6).Value = 217107
$ws.Cells.Item(323, 7).Value = 3109
$ws.Cells.Item(324, 6).Value = 249785
$ws.Cells.Item(324, 7).Value = 2857
$ws.Cells.Item(325, 6).Value = 774465
$ws.Cells.Item(326, 6).Value = 418125
$ws.Cells.Item(326, 7).Value = 3808
$ws.Cells.Item(328, 6).Value = 180924
$ws.Cells.Item(328, 7).Value = 2674
$ws.Cells.Item(329, 6).Value = 73257
$ws.Cells.Item(329, 7).Value = 1726
$ws.Cells.Item(334, 6).Value = 192976
$ws.Cells.Item(334, 7).Value = 3507
$ws.Cells.Item(335, 6).Value = 150063
$ws.Cells.Item(335, 7).Value = 3746
$ws.Cells.Item(337, 6).Value = 103542
$ws.Cells.Item(338, 6).Value = 221134
$ws.Cells.Item(341, 6).Value = 283668
$ws.Cells.Item(341, 7).Value = 3613
$ws.Cells.Item(342, 6).Value = 178588
$ws.Cells.Item(342, 7).Value = 3037
$ws.Cells.Item(349, 6).Value = 159407
$ws.Cells.Item(349, 7).Value = 2759
$ws.Cells.Item(350, 6).Value = 127028
$ws.Cells.Item(350, 7).Value = 2795
$ws.Cells.Item(351, 6).Value = 150738
$ws.Cells.Item(351, 7).Value = 2810
$ws.Cells.Item(352, 6).Value = 307397
$ws.Cells.Item(352, 7).Value = 3549
$ws.Cells.Item(355, 6).Value = 221975
$ws.Cells.Item(355, 7).Value = 3437
$ws.Cells.Item(356, 6).Value = 159997
$ws.Cells.Item(356, 7).Value = 2883
$ws.Cells.Item(357, 6).Value = 138221
$ws.Cells.Item(357, 7).Value = 3013
$ws.Cells.Item(358, 6).Value = 158805
$ws.Cells.Item(358, 7).Value = 2612
$ws.Cells.Item(359, 6).Value = 321351
$ws.Cells.Item(359, 7).Value = 3335
$ws.Cells.Item(360, 6).Value = 749755
$ws.Cells.Item(360, 7).Value = 5140
$ws.Cells.Item(367, 6).Value = 766547
$ws.Cells.Item(367, 7).Value = 3922
$ws.Cells.Item(371, 6).Value = 160185
$ws.Cells.Item(371, 7).Value = 1965
$ws.Cells.Item(384, 6).Value = 171574
$ws.Cells.Item(390, 6).Value = 219781

# Append new row 427 for 2021-05-06 (44321)
$ws.Range("A427").Value = 44321
$ws.Range("A427").NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(427, 2).Value = 384660
$ws.Cells.Item(427, 3).Value = 6268
$ws.Cells.Item(427, 4).Value = 343
$ws.Cells.Item(427, 5).Value = 11920
$ws.Cells.Item(427, 6).Value = 73202
$ws.Cells.Item(427, 7).Value = 295
